# Updates the PEBCOM worksheet with two new claim rows (55 and 56),
# matching the "Actualización automática del index.html y archivo Excel"
# commit: two rows appended to the bottom of the data table.
#
# The source data (row 1) is stored as plain text for every column except
# Coordenada_X / Coordenada_Y (M, N) which are numeric. Excel's COM layer
# auto-detects "numeric looking" strings (e.g. "6269", "4", "-490",
# "6/26/2025") and silently coerces them to numbers/dates when assigned
# through .Value. To keep them as genuine text values (as the rest of the
# sheet already does) we briefly force the cell to Text format before the
# assignment, then clear the format again so no stray style is left behind
# on the new cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($r, $c, [string]$text) {
    # Re-fetch the Range/Cell from $ws on every call (passing the COM object
    # itself through a function parameter does not keep a live binding in
    # this host), force Text format so numeric-looking strings ("6269",
    # "4", "6/26/2025", ...) are not auto-coerced into numbers/dates, then
    # clear the format again so the new cell ends up on the default style
    # (no stray "s" attribute), matching the rest of the sheet.
    $cell = $ws.Cells.Item($r, $c)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

$row55 = @("6269", "6/26/2025", "VELEZ SARSFIELD AV. 855", "4", "807789683", "PEBCOM", "Pendiente", "Columna inclinada", "1", "Cambio", "Sin equipos", "Pasante")
$row56 = @("-490", "6/25/2025", "Luzuriaga 1273", "4", "807789692", "PEBCOM", "Pendiente", "Picada", "1", "Cambio", "Sin equipos", "Pasante")

for ($col = 1; $col -le 12; $col++) {
    Set-TextValue 55 $col $row55[$col - 1]
    Set-TextValue 56 $col $row56[$col - 1]
}

$ws.Cells.Item(55, 13).Value = -58.389598
$ws.Cells.Item(55, 14).Value = -34.645174

$ws.Cells.Item(56, 13).Value = -58.387569
$ws.Cells.Item(56, 14).Value = -34.649344

Write-Host "Added rows 55-56 to PEBCOM sheet"
